$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D23").Value = "안녕하세요! 이번에 끝난 캐글 대회 RANZCR CLiP - Catheter and Line Position Challenge 에서 11위/15"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2760"

$ws.Range("D42").Value = "Python 인스타그램 크롤링 프로젝트 구조"
$ws.Range("E42").Value = "https://kjk92.tistory.com/64"

$ws.Range("D51").Value = "[git] git restore 사용법 정리"
$ws.Range("E51").Value = "https://bskyvision.com/1146"
